$d = $word.ActiveDocument

# Locate the two target paragraphs after the last table:
#  - the empty paragraph followed by "Schema logique :" (replaced with "Fonction :" + 4 equations)
#  - the paragraph holding the inline picture + bookmark (replaced with anchored picture + text, and bookmark split into its own paragraph)

$emptyPara = $null
$schemaPara = $null
$picturePara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.Contains("logique")) {
        $schemaPara = $pp
        $emptyPara = $d.Paragraphs.Item($i - 1)
        $picturePara = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($schemaPara -eq $null) {
    throw "Could not locate the 'logique' paragraph"
}

# --- Block 1: replace [empty paragraph .. "Schema logique :" paragraph] with "Fonction :" + 4 equations ---
$startPos = $emptyPara.Range.Start
$endPos = $schemaPara.Range.End
$block1Range = $d.Range($startPos, $endPos)

$block1Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="wp14"><w:body><w:p>
  <w:pPr>
    <w:rPr>
      <w:szCs w:val="26"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:szCs w:val="26"/>
    </w:rPr>
    <w:t>Fonction :</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:noProof/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
  <m:oMathPara>
    <m:oMath>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>Droite avance=</m:t>
      </m:r>
      <m:d>
        <m:dPr>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:dPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Droite+Avant</m:t>
          </m:r>
        </m:e>
      </m:d>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>.</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Gauche</m:t>
          </m:r>
        </m:e>
      </m:acc>
    </m:oMath>
  </m:oMathPara>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
      <w:noProof/>
      <w:szCs w:val="26"/>
    </w:rPr>
  </w:pPr>
  <m:oMathPara>
    <m:oMath>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>Gauche avance=</m:t>
      </m:r>
      <m:d>
        <m:dPr>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:dPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Gauche+Avant</m:t>
          </m:r>
        </m:e>
      </m:d>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>.</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Droite</m:t>
          </m:r>
        </m:e>
      </m:acc>
    </m:oMath>
  </m:oMathPara>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
      <w:noProof/>
      <w:szCs w:val="26"/>
    </w:rPr>
  </w:pPr>
  <m:oMathPara>
    <m:oMath>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>Gauche/Droite (Random) =</m:t>
      </m:r>
      <m:d>
        <m:dPr>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:dPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Droite.Gauche</m:t>
          </m:r>
        </m:e>
      </m:d>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>+(</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Droite</m:t>
          </m:r>
        </m:e>
      </m:acc>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>.</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Gauche</m:t>
          </m:r>
        </m:e>
      </m:acc>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>.</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Avant</m:t>
          </m:r>
        </m:e>
      </m:acc>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>.</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Arrière</m:t>
          </m:r>
        </m:e>
      </m:acc>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>)</m:t>
      </m:r>
    </m:oMath>
  </m:oMathPara>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:noProof/>
      <w:sz w:val="18"/>
    </w:rPr>
  </w:pPr>
  <m:oMathPara>
    <m:oMath>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>Droite recule=</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Droite</m:t>
          </m:r>
        </m:e>
      </m:acc>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>.</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Gauche</m:t>
          </m:r>
        </m:e>
      </m:acc>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>.</m:t>
      </m:r>
      <m:acc>
        <m:accPr>
          <m:chr m:val="̅"/>
          <m:ctrlPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:i/>
              <w:szCs w:val="26"/>
            </w:rPr>
          </m:ctrlPr>
        </m:accPr>
        <m:e>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
              <w:szCs w:val="26"/>
            </w:rPr>
            <m:t>Avant</m:t>
          </m:r>
        </m:e>
      </m:acc>
      <m:r>
        <w:rPr>
          <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <m:t>.Arrière</m:t>
      </m:r>
    </m:oMath>
  </m:oMathPara>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$block1Range.InsertXML($block1Xml)

# --- Block 2: replace the picture paragraph with anchored picture + "Schema logique :" run + new bookmark paragraph ---
# Re-resolve the picture paragraph since the document structure shifted above.
$picturePara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.InlineShapes.Count -gt 0) {
        $picturePara2 = $pp
    }
}
if ($picturePara2 -eq $null) {
    throw "Could not locate the picture paragraph"
}

$block2Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="wp14"><w:body><w:p>
  <w:r>
    <w:rPr>
      <w:noProof/>
    </w:rPr>
    <w:drawing>
      <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251658240" behindDoc="1" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3D5F032F">
        <wp:simplePos x="0" y="0"/>
        <wp:positionH relativeFrom="margin">
          <wp:align>center</wp:align>
        </wp:positionH>
        <wp:positionV relativeFrom="paragraph">
          <wp:posOffset>263525</wp:posOffset>
        </wp:positionV>
        <wp:extent cx="4395470" cy="2571750"/>
        <wp:effectExtent l="0" t="0" r="5080" b="0"/>
        <wp:wrapTight wrapText="bothSides">
          <wp:wrapPolygon edited="0">
            <wp:start x="0" y="0"/>
            <wp:lineTo x="0" y="21440"/>
            <wp:lineTo x="21531" y="21440"/>
            <wp:lineTo x="21531" y="0"/>
            <wp:lineTo x="0" y="0"/>
          </wp:wrapPolygon>
        </wp:wrapTight>
        <wp:docPr id="1" name="Image 1"/>
        <wp:cNvGraphicFramePr>
          <a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/>
        </wp:cNvGraphicFramePr>
        <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
          <a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">
            <pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">
              <pic:nvPicPr>
                <pic:cNvPr id="1" name=""/>
                <pic:cNvPicPr/>
              </pic:nvPicPr>
              <pic:blipFill rotWithShape="1">
                <a:blip r:embed="rId4">
                  <a:extLst>
                    <a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}">
                      <a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/>
                    </a:ext>
                  </a:extLst>
                </a:blip>
                <a:srcRect l="50661" t="33157" r="1191" b="16755"/>
                <a:stretch/>
              </pic:blipFill>
              <pic:spPr bwMode="auto">
                <a:xfrm>
                  <a:off x="0" y="0"/>
                  <a:ext cx="4395470" cy="2571750"/>
                </a:xfrm>
                <a:prstGeom prst="rect">
                  <a:avLst/>
                </a:prstGeom>
                <a:ln>
                  <a:noFill/>
                </a:ln>
                <a:extLst>
                  <a:ext uri="{53640926-AAD7-44D8-BBD7-CCE9431645EC}">
                    <a14:shadowObscured xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"/>
                  </a:ext>
                </a:extLst>
              </pic:spPr>
            </pic:pic>
          </a:graphicData>
        </a:graphic>
        <wp14:sizeRelH relativeFrom="margin">
          <wp14:pctWidth>0</wp14:pctWidth>
        </wp14:sizeRelH>
        <wp14:sizeRelV relativeFrom="margin">
          <wp14:pctHeight>0</wp14:pctHeight>
        </wp14:sizeRelV>
      </wp:anchor>
    </w:drawing>
  </w:r>
  <w:r>
    <w:t>Schéma logique :</w:t>
  </w:r>
</w:p>
<w:p>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$picturePara2.Range.InsertXML($block2Xml)

Write-Output "Done. Paragraph count:"
Write-Output $d.Paragraphs.Count
